$d = $word.ActiveDocument

# 1) Drop the leading "The " from the four law headings (but keep "Introduction" etc. untouched)
$d.Content.Find.Execute("The Four Laws of Thermodynamics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Four Laws of Thermodynamics", 2) | Out-Null

$d.Content.Find.Execute("The Zeroth Law of Thermodynamics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Zeroth Law of Thermodynamics", 2) | Out-Null

# The "First Law of Thermodynamics" heading is split across three runs with a
# _GoBack bookmark in the middle; replacing the whole visible phrase collapses
# it into a single run and removes that now-orphaned bookmark.
$d.Content.Find.Execute("The First Law of Thermodynamics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "First Law of Thermodynamics", 2) | Out-Null

$d.Content.Find.Execute("The Second Law of Thermodynamics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Second Law of Thermodynamics", 2) | Out-Null

$d.Content.Find.Execute("The Third Law of Thermodynamics", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Third Law of Thermodynamics", 2) | Out-Null

# 2) Re-create the _GoBack bookmark (Word always leaves one at the last edit
# location) right after the "Maxwell Relations" entry, i.e. at the end of
# that paragraph's text, before the paragraph mark.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $d.Content
$target.Find.Execute("Maxwell Relations") | Out-Null
$target.Collapse(0)

# Inserting a temporary marker character lets us anchor a zero-length
# bookmark precisely at the end of the text (COM mis-resolves a bookmark
# added directly on an already-collapsed range at this position), then we
# remove the marker again, leaving the bookmark collapsed in place.
$target.InsertAfter("X")
$target.Bookmarks.Add("_GoBack") | Out-Null
$target.Text = ""
